$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.352274417877197
$ws.Range("B1").Value = 1.40834367275238
$ws.Range("C1").Value = 3.91486668586731
$ws.Range("D1").Value = 3.166523456573486
$ws.Range("E1").Value = 1.056880950927734
